$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): G2 becomes text "不可售"; several F-column counters increment
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = "不可售"
$ws1.Range("F3").Value = 377
$ws1.Range("F4").Value = 4830
$ws1.Range("F5").Value = 12
$ws1.Range("F7").Value = 484

# Sheet "全部类型" (sheet4): same logical rows, but shifted down by the extra
# row present in this sheet (row 5 here corresponds to row 4 in 展览, etc.)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F3").Value = 377
$ws4.Range("F4").Value = 4830
$ws4.Range("F6").Value = 12
$ws4.Range("F9").Value = 484
